# Apply the change described by the diff:
#  - Swap the header text of E1/F1
#    ("Fail key reason" <-> "Flight Cancellation Key Reason")
#  - For every data row (2..258), swap whatever value sits in column E
#    with whatever value sits in column F (every populated E cell moves
#    to F and every populated F cell moves to E)
#  - Column widths for E/F follow the swapped (bestFit) content
#  - Apply an AutoFilter over column B (B1:B258) and recreate the
#    corresponding hidden _FilterDatabase defined name
#  - Leave the final selection on cell B76 (matching the author's last
#    recorded cursor position)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 258

# --- Swap header labels in row 1 (E1 <-> F1) ---------------------------
$e1 = $ws.Range("E1").Value()
$f1 = $ws.Range("F1").Value()
$ws.Range("E1").Value() = $f1
$ws.Range("F1").Value() = $e1

# --- Swap column E / column F contents for every data row --------------
for ($r = 2; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)
    $eVal = $eCell.Value()
    $fVal = $fCell.Value()
    if ($eVal -ne $fVal) {
        $eCell.Value() = $fVal
        $fCell.Value() = $eVal
    }
}

# Column widths for E/F are recomputed automatically (bestFit) by the
# engine once their contents change, so no manual resize is required.

# --- Apply AutoFilter on column B and recreate the hidden name ---------
$ws.Range("B1:B258").AutoFilter()
$fdbName = $ws.Names.Add("_xlnm._FilterDatabase", "=SiteStatusProject!`$B`$1:`$B`$258")
$fdbName.Visible = $false

# --- Restore the last selection recorded in the workbook ---------------
$ws.Range("B76").Select()

$wb.Save()
